$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Record", "Balanço Geral", "Social", "2025-04-03T12:25", "Negativo", "Beneficiários do CadÚnico dormem na fila para tentar atualizar o cadastro. Repórter *ao vivo*. Hoje, segundo dia do mutirão, foi mais organizado do que ontem. Na fila, funcionários colocavam pulseira para identificar as prioridades. Apesar de estar melhor do que ontem, muitas pessoas ainda precisam dormir na fila para conseguir atendimento. Depoimento de beneficiária que mora em Travessão. Para se candidatar a um beneficio, tem que ir aos Cras. Mutirão é para desafogar os Cras. Entrevista com o assessor da Secretaria, Ruan Barros. "),
    @("Record", "Balanço Geral", "Defesa Civil", "2025-04-03T12:43", "Neutro", "Previsão de chuvas intensas para todo Estado do Rio de hoje até domingo. Repórter *ao vivo*. Vídeo com depoimento do climatologista Carlos Augusto Souto.  "),
    @("Record", "Balanço Geral", "Saúde", "2025-04-03T12:04", "Positivo", "Começou hoje a atualização cadastral de pacientes assistidos pelo SAD. Repórter *ao vivo*. Prazo vai até 5 de abril. Serão duas etapas: documentação e visita técnica domiciliar. Objetivo é ver a complexidade de cada paciente. As pessoas devem levar os documentos à sede da Secretaria de Saúde."),
    @("Record", "Balanço Geral", "Meio Ambiente", "2025-04-03T13:20", "Neutro", "Vazamento de esgoto gera mau cheiro e incomoda moradores do bairro da Penha. É na Travessa Rossini Quintanilha. Vazamento está se estendendo pela via causando mau cheiro. Exibido vídeo enviado por morador, que alegou que o *problema teve início desde que foi feito o Bairro Legal da Penha*. Problema dura seis anos. Concessionária foi procurada, mas não respondeu. ")
)

$startRow = 28
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
